$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted columns (B, C, D) keep their original "Text" number format so that
# numeric-looking strings (e.g. "313.39") are not silently converted to numbers by Excel.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.697.75"
$ws.Range("E2").Value = "  -0.36%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.543.76"
$ws.Range("E3").Value = "  +0.56%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.39"
$ws.Range("E5").Value = "  -1.41%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.26"
$ws.Range("E6").Value = "  +4.68%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.570"
$ws.Range("E7").Value = "  -0.67%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").Value = "  -0.93%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.17"
$ws.Range("E10").Value = "  +0.91%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0806"
$ws.Range("E11").Value = "  -1.22%  "

# Row 12
$ws.Range("E12").Value = "  -1.68%  "

# Row 13
$ws.Range("E13").Value = "  -0.04%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.937.24"
$ws.Range("E14").Value = "  +0.68%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.82"
$ws.Range("E15").Value = "  +5.02%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.584.83"
$ws.Range("E16").Value = "  +0.90%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.838"
$ws.Range("E17").Value = "  -1.26%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.722.80"
$ws.Range("E18").Value = "  -0.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.77"
$ws.Range("E19").Value = "  -1.12%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.36"
$ws.Range("E20").Value = "  -1.63%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0955"
$ws.Range("E21").Value = "  -1.00%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.10"
$ws.Range("E22").Value = "  -0.71%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "246.33"
$ws.Range("E23").Value = "  -2.58%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.94"
$ws.Range("E24").Value = "  -0.85%  "

# Row 25
$ws.Range("E25").Value = "  +0.39%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.48"
$ws.Range("E26").Value = "  +0.37%  "

# Row 27
$ws.Range("E27").Value = "  -0.03%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40.93"
$ws.Range("E28").Value = "  -1.04%  "

# Row 29
$ws.Range("E29").Value = "  -2.20%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.14"
$ws.Range("E30").Value = "  -2.68%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.69"
$ws.Range("E31").Value = "  -0.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.73"
$ws.Range("E32").Value = "  -2.49%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0804"
$ws.Range("E33").Value = "  +1.76%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.08"
$ws.Range("E34").Value = "  -3.16%  "

# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "ApeXProtocol"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.64"
$ws.Range("E35").Value = "  +7.90%  "

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.26"
$ws.Range("E36").Value = "  -3.03%  "

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.64"
$ws.Range("E37").Value = "  -2.85%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.31"
$ws.Range("E38").Value = "  -5.52%  "

# Row 39
$ws.Range("E39").Value = "  -1.27%  "

# Row 40
$ws.Range("E40").Value = "  -0.89%  "

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.19"
$ws.Range("E41").Value = "  +10.16%  "

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.33"
$ws.Range("E42").Value = "  +1.96%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0299"
$ws.Range("E44").Value = "  -1.85%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.28"
$ws.Range("E45").Value = "  +0.09%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.980.77"
$ws.Range("E46").Value = "  -0.88%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.01"
$ws.Range("E47").Value = "  -1.24%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.793.43"
$ws.Range("E48").Value = "  +0.82%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "81.39"
$ws.Range("E49").Value = "  -3.58%  "

# Row 50
$ws.Range("E50").Value = "  +0.50%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.53"
$ws.Range("E51").Value = "  -2.01%  "
